$d = $word.ActiveDocument

$RED = 255
$AUTO = -16777216

function Get-ParaRange($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    return $d.Range($p.Range.Start, $p.Range.End)
}

function Replace-InPara($paraIndex, $findText, $replaceText) {
    $rng = Get-ParaRange $paraIndex
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1)
    if (-not $found) {
        Write-Host "WARNING: not found '$findText' in paragraph $paraIndex"
    }
    return $found
}

# Find `fullText` (unique) in paragraph, then color the sub-span [0, splitAt) auto
# and [splitAt, end) red (relative offsets into fullText, measured in characters).
function Split-Color($paraIndex, $fullText, $splitAt, $autoColor, $redColor) {
    $rng = Get-ParaRange $paraIndex
    $found = $rng.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: split-find not found '$fullText' in paragraph $paraIndex"
        return
    }
    $start = $rng.Start
    $end = $rng.End
    if ($splitAt -gt 0) {
        $r1 = $d.Range($start, $start + $splitAt)
        $r1.Font.Color = $autoColor
    }
    $r2 = $d.Range($start + $splitAt, $end)
    $r2.Font.Color = $redColor
}

# Color the whole found occurrence of `fullText` a single color.
function Color-Whole($paraIndex, $fullText, $colorVal) {
    $rng = Get-ParaRange $paraIndex
    $found = $rng.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: color-find not found '$fullText' in paragraph $paraIndex"
        return
    }
    $rng.Font.Color = $colorVal
}

function Append-ToPara($paraIndex, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $e = $p.Range.End
    $textEnd = $e - 1
    $insPoint = $d.Range($textEnd, $textEnd)
    $insPoint.InsertAfter($text)
}

# ---- Paragraph 38: "Visitante" -> "Visita" ----
Replace-InPara 38 "Visitante" "Visita"

# ---- Paragraph 39: "Novo Visitante" -> "Nova Visita"; "Tela 0201" -> "Tela 1117" (split colors) ----
Replace-InPara 39 "Novo Visitante" "Nova Visita"
Replace-InPara 39 "Tela 0201" "Tela 1117"
Split-Color 39 "Tela 1117" 4 $AUTO $RED

# ---- Paragraph 40: "Tela 0301" -> "Tela 1118" (split colors: "Tela " auto, "1118" red) ----
Replace-InPara 40 "Tela 0301" "Tela 1118"
Split-Color 40 "Tela 1118" 5 $AUTO $RED

# ---- Paragraph 50: append "(Tela 1119)" after "passo 3." ----
Append-ToPara 50 "(Tela 1119)"
Color-Whole 50 "Tela 1119" $RED

# ---- Paragraph 59: standalone "0201" -> "1117" ----
Replace-InPara 59 "0201" "1117"

# ---- Paragraph 60: standalone "0501" -> "1120" ----
Replace-InPara 60 "0501" "1120"

# ---- Paragraph 74: "Tela 0501" -> "Tela 1120"; "Tela 0901" -> "Tela 1121" ----
Replace-InPara 74 "Tela 0501" "Tela 1120"
Split-Color 74 "Tela 1120" 5 $AUTO $RED
Replace-InPara 74 "Tela 0901" "Tela 1121"
Split-Color 74 "Tela 1121" 5 $AUTO $RED

# ---- Paragraph 84: gender fix + "O" -> "A"; "Tela 1001" -> "Tela 1122" ----
Replace-InPara 84 "excluído. Sistema exibe mensagem “O" "excluída. Sistema exibe mensagem “A"
Replace-InPara 84 "Tela 1001" "Tela 1122"
Split-Color 84 "Tela 1122" 7 $AUTO $RED

# ---- Paragraph 93: standalone "0201" -> "1117" ----
Replace-InPara 93 "0201" "1117"

# ---- Paragraph 94: "Tela 0401" -> "Tela 1123" ----
Replace-InPara 94 "Tela 0401" "Tela 1123"
Split-Color 94 "Tela 1123" 5 $AUTO $RED

# ---- Paragraph 100: append " (Tela 1117)." at end ----
Append-ToPara 100 " (Tela 1117)."
Color-Whole 100 "Tela 1117" $RED

# ---- Paragraph 104: append "(Tela 1119)" after "passo 2." ----
Append-ToPara 104 "(Tela 1119)"
Color-Whole 104 "Tela 1119" $RED

Write-Host "DONE"
